# Build the "ValidLogin" / "Invalidlogin" login-demo workbook.
$wb = $excel.ActiveWorkbook

# --- Rename the first two sheets ---------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Invalidlogin"

# --- ValidLogin sheet: Username/Password table of valid accounts -------
$ws1.Range("A1").Value = "Username"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"
$ws1.Range("A3").Value = "trainee"
$ws1.Range("B3").Value = "trainee"

$ws1.Columns.Item(1).ColumnWidth = 10

# --- Invalidlogin sheet: Username/Password example of an invalid login -
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xyz"

# --- Selections matching the authored workbook --------------------------
$ws1.Range("E15").Select()

$ws2.Activate()
$ws2.Range("H18").Select()
